# Auto update Excel log
#
# Appends the next batch of sensor-log readings (timestamped 2026-02-06,
# just after the previous last row on each sheet) to the PIR, Humidity and
# Temperature sheets. Columns are: Date, Timestamp, Hour, Location, Value,
# Status.
#
# Every value is written with a leading single-quote so Excel stores it as
# literal text, matching the existing rows (plain text, "General" format)
# instead of auto-converting date-looking strings (2026-02-06) to real
# dates or percent-looking strings (69.4%) to numeric percentages.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PIR sheet: append rows 297-310
# ---------------------------------------------------------------------------
$pirSheet = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-02-06", "10:04:50", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:04:52", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:04:55", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:05:00", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:05:02", "10:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "10:05:10", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:05:15", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:05:20", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:05:25", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:05:28", "10:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "10:05:37", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:05:38", "10:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "10:05:47", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:05:48", "10:00", "Bathroom", "Motion Detected", "Active")
)
$r = 297
foreach ($row in $pirRows) {
    for ($c = 1; $c -le 6; $c++) {
        $pirSheet.Cells.Item($r, $c).Value = "'" + $row[$c - 1]
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 192-204
# ---------------------------------------------------------------------------
$humiditySheet = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-02-06", "10:04:51", "10:00", "Bathroom", "69.4%", "Active"),
    @("2026-02-06", "10:04:53", "10:00", "Bathroom", "69.3%", "Active"),
    @("2026-02-06", "10:04:56", "10:00", "Bathroom", "69.4%", "Active"),
    @("2026-02-06", "10:05:01", "10:00", "Bathroom", "69.4%", "Active"),
    @("2026-02-06", "10:05:06", "10:00", "Bathroom", "69.5%", "Active"),
    @("2026-02-06", "10:05:11", "10:00", "Bathroom", "69.5%", "Active"),
    @("2026-02-06", "10:05:16", "10:00", "Bathroom", "69.4%", "Active"),
    @("2026-02-06", "10:05:21", "10:00", "Bathroom", "69.7%", "Active"),
    @("2026-02-06", "10:05:26", "10:00", "Bathroom", "69.8%", "Active"),
    @("2026-02-06", "10:05:31", "10:00", "Bathroom", "69.9%", "Active"),
    @("2026-02-06", "10:05:36", "10:00", "Bathroom", "70.1%", "Active"),
    @("2026-02-06", "10:05:41", "10:00", "Bathroom", "70.3%", "Active"),
    @("2026-02-06", "10:05:46", "10:00", "Bathroom", "70.2%", "Active")
)
$r = 192
foreach ($row in $humidityRows) {
    for ($c = 1; $c -le 6; $c++) {
        $humiditySheet.Cells.Item($r, $c).Value = "'" + $row[$c - 1]
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Temperature sheet: append rows 192-204
# ---------------------------------------------------------------------------
$temperatureSheet = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-02-06", "10:04:52", "10:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "10:04:54", "10:00", "Bathroom", "27.6C", "Active"),
    @("2026-02-06", "10:04:56", "10:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "10:05:01", "10:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "10:05:06", "10:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "10:05:11", "10:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "10:05:16", "10:00", "Bathroom", "27.7C", "Active"),
    @("2026-02-06", "10:05:21", "10:00", "Bathroom", "27.8C", "Active"),
    @("2026-02-06", "10:05:26", "10:00", "Bathroom", "27.8C", "Active"),
    @("2026-02-06", "10:05:31", "10:00", "Bathroom", "27.8C", "Active"),
    @("2026-02-06", "10:05:36", "10:00", "Bathroom", "27.8C", "Active"),
    @("2026-02-06", "10:05:41", "10:00", "Bathroom", "27.8C", "Active"),
    @("2026-02-06", "10:05:46", "10:00", "Bathroom", "27.9C", "Active")
)
$r = 192
foreach ($row in $temperatureRows) {
    for ($c = 1; $c -le 6; $c++) {
        $temperatureSheet.Cells.Item($r, $c).Value = "'" + $row[$c - 1]
    }
    $r = $r + 1
}
